# Update the "取得日時" (acquisition timestamp) column for all existing
# data rows on the "ランサーズ" sheet to reflect the new run time.
# Diff: 2026-01-21 12:44:19  ->  2026-01-21 12:58:06  (rows 2-8, column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-21 12:58:06"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
